$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'262.91"
$ws.Range("D3").Value = "'22.75"
$ws.Range("D4").Value = "'6.199"
$ws.Range("D5").Value = "'0.06127"
$ws.Range("D6").Value = "'3.511"
$ws.Range("D7").Value = "'6.703"
$ws.Range("D8").Value = "'1.358"
$ws.Range("D9").Value = "'0.7994"
$ws.Range("D10").Value = "'0.1573"
$ws.Range("D11").Value = "'0.08143"
$ws.Range("D12").Value = "'0.03323"
$ws.Range("D13").Value = "'0.03155"
$ws.Range("D14").Value = "'0.09259"
$ws.Range("D15").Value = "'3.934"
$ws.Range("D16").Value = "'0.001688"
$ws.Range("D17").Value = "'0.04828"
$ws.Range("D18").Value = "'0.0006210"
$ws.Range("D19").Value = "'0.006195"
$ws.Range("D20").Value = "'0.001100"
$ws.Range("D21").Value = "'0.003195"
$ws.Range("D24").Value = "'2.280"
$ws.Range("D25").Value = "'0.3374"
$ws.Range("D26").Value = "'0.1251"
$ws.Range("D27").Value = "'0.0002682"
$ws.Range("D40").Value = "'0.04598"
$ws.Range("D41").Value = "'0.007263"
$ws.Range("D42").Value = "'0.003902"
$ws.Range("D44").Value = "'0.01087"
$ws.Range("D45").Value = "'0.002970"
$ws.Range("D46").Value = "'0.00005996"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D48").Value = "'0.7000"
$ws.Range("D49").Value = "'0.04957"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D51").Value = "'0.01010"
